# --- Setup -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$param = $wb.Worksheets.Item("Param")

# workbook-level: fix the absPath + activeTab will follow once new sheet is active

# Add the new worksheet right after "Param"
$ws = $wb.Worksheets.Add($null, $param)
$ws.Name = "Parameters of interest"

# Colors (Excel COM color longs are 0xBBGGRR, i.e. R + G*256 + B*65536)
$colPink   = 13289973   # F5C9CA
$colWhite  = 16777215   # FFFFFF
$colBorder = 2104739    # A31D20
$colBlack  = 0          # 000000

# Border line-style / weight constants
$xlContinuous = 1
$xlMedium = -4138

# Alignment constants
$xlLeft = -4131
$xlCenter = -4108

function Style-ALabel($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Size = 20
    $rng.Font.Color = $colBlack
    $rng.Font.Name = "Microsoft New Tai Lue"
    $rng.Interior.Color = $colPink
    $rng.Borders.Item(7).LineStyle = $xlContinuous
    $rng.Borders.Item(7).Weight = $xlMedium
    $rng.Borders.Item(7).Color = $colBorder
    $rng.Borders.Item(8).LineStyle = $xlContinuous
    $rng.Borders.Item(8).Weight = $xlMedium
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = $xlContinuous
    $rng.Borders.Item(9).Weight = $xlMedium
    $rng.Borders.Item(9).Color = $colBorder
    $rng.HorizontalAlignment = $xlLeft
    $rng.WrapText = $true
    $rng.RowHeight = 31
}

function Style-BValue($rng, [bool]$bold) {
    $rng.Font.Bold = $bold
    $rng.Font.Size = 20
    $rng.Font.Color = $colBlack
    $rng.Font.Name = "Microsoft New Tai Lue"
    $rng.Interior.Color = $colWhite
    $rng.Borders.Item(10).LineStyle = $xlContinuous
    $rng.Borders.Item(10).Weight = $xlMedium
    $rng.Borders.Item(10).Color = $colBorder
    $rng.Borders.Item(8).LineStyle = $xlContinuous
    $rng.Borders.Item(8).Weight = $xlMedium
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = $xlContinuous
    $rng.Borders.Item(9).Weight = $xlMedium
    $rng.Borders.Item(9).Color = $colBorder
    $rng.HorizontalAlignment = $xlCenter
    $rng.WrapText = $true
    $rng.RowHeight = 31
}

function Format-Run($cellRange, [int]$start, [int]$len, [bool]$bold) {
    $c = $cellRange.Characters($start, $len)
    $c.Font.Bold = $bold
    $c.Font.Size = 20
    $c.Font.Color = $colBlack
    $c.Font.Name = "Microsoft New Tai Lue"
}

# --- Row 1: PathToFasta / full_yeast ---------------------------------------
$ws.Range("A1").Value = "PathToFasta"
Style-ALabel $ws.Range("A1")
$ws.Range("B1").Value = "full_yeast"
Style-BValue $ws.Range("B1") $true

# --- Row 2: NumReps / 2,3,4,5,6 (3 bold) ------------------------------------
$ws.Range("A2").Value = "NumReps"
Style-ALabel $ws.Range("A2")
$ws.Range("B2").Value = "2,3,4,5,6"
Style-BValue $ws.Range("B2") $false
Format-Run $ws.Range("B2") 3 1 $true
Format-Run $ws.Range("B2") 4 6 $false

# --- Row 3: NumCond / 2 ------------------------------------------------------
$ws.Range("A3").Value = "NumCond"
Style-ALabel $ws.Range("A3")
$ws.Range("B3").Value = 2
Style-BValue $ws.Range("B3") $true

# --- Row 4: QuantNoise / 0.01,0.125,0.25 (0.25 bold) ------------------------
$ws.Range("A4").Value = "QuantNoise"
Style-ALabel $ws.Range("A4")
$ws.Range("B4").Value = "0.01,0.125,0.25"
Style-BValue $ws.Range("B4") $false
Format-Run $ws.Range("B4") 12 4 $true

# --- Row 5: DiffRegFrac / 0.1 -----------------------------------------------
$ws.Range("A5").Value = "DiffRegFrac"
Style-ALabel $ws.Range("A5")
$ws.Range("B5").Value = 0.1
Style-BValue $ws.Range("B5") $true

# --- Row 6: DiffRegMax / 2 ---------------------------------------------------
$ws.Range("A6").Value = "DiffRegMax"
Style-ALabel $ws.Range("A6")
$ws.Range("B6").Value = 2
Style-BValue $ws.Range("B6") $true

# --- Row 7: Enzyme / trypsin, trypsin.strict (trypsin bold) -----------------
$ws.Range("A7").Value = "Enzyme"
Style-ALabel $ws.Range("A7")
$ws.Range("B7").Value = "trypsin, trypsin.strict"
Style-BValue $ws.Range("B7") $true
Format-Run $ws.Range("B7") 8 16 $false

# --- Row 8: PropMissedCleavages / 0.01,0.1,0.2 (0.01 bold) ------------------
$ws.Range("A8").Value = "PropMissedCleavages"
Style-ALabel $ws.Range("A8")
$ws.Range("B8").Value = "0.01,0.1,0.2"
Style-BValue $ws.Range("B8") $true
Format-Run $ws.Range("B8") 5 8 $false

# --- Row 9: PercDetectedPep / 0.1,0.25,0.5 (0.5 bold) -----------------------
$ws.Range("A9").Value = "PercDetectedPep"
Style-ALabel $ws.Range("A9")
$ws.Range("B9").Value = "0.1,0.25,0.5"
Style-BValue $ws.Range("B9") $false
Format-Run $ws.Range("B9") 10 3 $true

# --- Row 10: PercDetectedVal / (same rich string as row 9) ------------------
$ws.Range("A10").Value = "PercDetectedVal"
Style-ALabel $ws.Range("A10")
$ws.Range("B9").Copy($ws.Range("B10"))

# --- Row 11: MSNoise / 0.01,0.125,0.25 (0.125 bold) -------------------------
$ws.Range("A11").Value = "MSNoise"
Style-ALabel $ws.Range("A11")
$ws.Range("B11").Value = "0.01,0.125,0.25"
Style-BValue $ws.Range("B11") $false
Format-Run $ws.Range("B11") 6 5 $true
Format-Run $ws.Range("B11") 11 5 $false

# --- Row 12: WrongIDs / 0.01,0.05 (0.01 bold) --------------------------------
$ws.Range("A12").Value = "WrongIDs"
Style-ALabel $ws.Range("A12")
$ws.Range("B12").Value = "0.01,0.05"
Style-BValue $ws.Range("B12") $true
Format-Run $ws.Range("B12") 5 5 $false

# --- Sheet-level view/selection tweaks --------------------------------------
$ws.Range("A14").Select()

# restore Param's own selection/tab state
$param.Range("A8").Select()
$param.Activate()

Write-Host "done"
